# ------------------------------------------------------------------
# Edit script: rewrite the "consciousness" essay into the "numbers"
# essay, rename the author/email, fix the TimesNewToman -> Times New
# Roman font typo, and append a trailing blank paragraph.
# ------------------------------------------------------------------

$d = $word.ActiveDocument
$vt = [char]11   # w:br (line break) marker inside Range.Text

function Replace-ParagraphText($paraIndex, $newText) {
    $para = $d.Paragraphs($paraIndex).Range
    $start = $para.Start
    $end = $para.End
    if ($end -gt $start) {
        $head = $d.Range($start, $start + 1)
        $head.Text = $newText
    } else {
        $para.InsertAfter($newText)
    }
    $newEnd = $start + $newText.Length
    $paraNow = $d.Paragraphs($paraIndex).Range
    if ($paraNow.End -gt $newEnd) {
        $remainder = $d.Range($newEnd, $paraNow.End)
        $remainder.Delete()
    }
}

# 1) Fix the misspelled font name everywhere (format-based, not text-based).
$fontFind = $d.Content.Find
$fontFind.ClearFormatting()
$fontFind.Replacement.ClearFormatting()
$fontFind.Font.Name = "TimesNewToman"
$fontFind.Replacement.Font.Name = "Times New Roman"
$fontFind.Execute("", $false, $false, $false, $false, $false, $true, 1, $true, "", 2) | Out-Null

# 2) Title.
$titleFind = $d.Content.Find
$titleFind.ClearFormatting()
$titleFind.Replacement.ClearFormatting()
$titleFind.Execute('Unraveling the Enigma of Consciousness', $true, $false, $false, $false, $false, $true, 1, $false, 'The Allure of Numbers: A Mathematical Journey through Patterns, Equations, and Problem Solving', 2) | Out-Null

# 3) Author name.
$nameFind = $d.Content.Find
$nameFind.ClearFormatting()
$nameFind.Replacement.ClearFormatting()
$nameFind.Execute('Eleanor Gamaliel', $true, $false, $false, $false, $false, $true, 1, $false, 'Eleanor Hayes', 2) | Out-Null

# 4) Email address (paragraph 3): "e" + "." + "gamaliel@research" + "." + "hub"
#    becomes "eleanorhayes789@protonmail" + "." + "com".
$emailPara = $d.Paragraphs(3).Range
$firstChar = $d.Range($emailPara.Start, $emailPara.Start + 1)
$firstChar.Text = 'eleanorhayes789@protonmail'

$domainRange = $d.Paragraphs(3).Range
$domainFind = $domainRange.Find
$domainFind.ClearFormatting()
$domainFind.Replacement.ClearFormatting()
$domainFind.Execute('gamaliel@research', $true, $false, $false, $false, $false, $true, 1, $false, 'com', 2) | Out-Null

$tailRange = $d.Paragraphs(3).Range
$tailFind = $tailRange.Find
$tailFind.ClearFormatting()
$tailFind.Replacement.ClearFormatting()
$tailFound = $tailFind.Execute('.hub')
if ($tailFound) {
    $tailRange.Delete()
}

# 5) Body paragraph (paragraph 5).
$bodyText = 'In the realm of human knowledge, mathematics stands as a beacon of precision, logic, and elegance. It is a language of patterns, equations, and problem solving that has captivated curious minds for centuries. Like a symphony of numbers, mathematics unveils the underlying order in our universe, guiding us through the complexities of modern life. Its beauty lies in its universality, transcending cultural boundaries and connecting individuals across time and space. As we embark on this mathematical journey, we will explore the intricate world of numbers, uncovering the remarkable power of patterns, the elegance of equations, and the satisfaction of solving challenging problems.' + $vt + '' + $vt + 'In the vast tapestry of mathematics, patterns emerge as fundamental building blocks. From the mesmerizing symmetry of snowflakes to the Fibonacci sequence found in nature''s spirals, patterns reveal a hidden order that governs the universe. By studying patterns, mathematicians seek to uncover the underlying principles that shape our world, unraveling the enigmatic puzzles that have intrigued humanity for generations. Through the exploration of patterns, we gain insights into the intricate workings of our surroundings and develop a deeper appreciation for the interconnectedness of all things.' + $vt + '' + $vt + 'Equations, the lifeblood of mathematics, provide a concise and powerful way to express complex relationships between quantities. They allow us to model real-world phenomena, from the trajectory of a projectile to the intricate dynamics of weather systems. Equations encode the essence of mathematical truths, revealing hidden connections and unlocking the secrets of the universe. By manipulating and solving equations, we gain insights into the fundamental forces that govern our world and develop the ability to predict and control various aspects of our lives.' + $vt + '' + $vt + 'Problem solving, the ultimate test of mathematical prowess, challenges us to apply our knowledge and skills to untangle complex and perplexing situations. Mathematical problems come in various forms, from brain-teasers and puzzles to real-world challenges that require creative and analytical thinking. Solving problems not only sharpens our minds and develops our problem-solving skills but also instills a sense of accomplishment and satisfaction. Through problem solving, we learn to persevere in the face of challenges, think critically, and approach problems from different angles, skills that are invaluable in all aspects of life.'
Replace-ParagraphText 5 $bodyText

# 6) Summary paragraph (paragraph 7).
$summaryText = 'Mathematics is a vast and captivating field of study that encompasses patterns, equations, and problem solving. Patterns reveal the underlying order in our universe, equations provide a concise and powerful way to express complex relationships, and problem solving challenges us to apply our knowledge and skills to untangle complex situations. Through mathematics, we gain insights into the fundamental forces that govern our world, develop valuable problem-solving skills, and cultivate a deeper appreciation for the elegance and beauty of the numerical realm. Mathematics is a journey that continues to inspire and intrigue, inviting us to explore the hidden depths of the universe and unlock the secrets it holds.'
Replace-ParagraphText 7 $summaryText

# 7) Append a trailing blank paragraph at the very end of the document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

Write-Host "done"
